$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New shared-string column headers (row 1), columns CD:CI ---
$ws.Range("CD1").Value2 = "sp.mnRng"
$ws.Range("CE1").Value2 = "gen.mnRng"
$ws.Range("CF1").Value2 = "sf.mnRng"
$ws.Range("CG1").Value2 = "sp.medRng"
$ws.Range("CH1").Value2 = "gen.medRng"
$ws.Range("CI1").Value2 = "sf.medRng"
$ws.Range("CD2").Value2 = 924.62818791946302
$ws.Range("CE2").Value2 = 1275.08
$ws.Range("CF2").Value2 = 1625.6
$ws.Range("CG2").Value2 = 914.4
$ws.Range("CH2").Value2 = 1295.4000000000001
$ws.Range("CI2").Value2 = 1524
$ws.Range("CD3").Value2 = 951.24301487603304
$ws.Range("CE3").Value2 = 1526.34563478261
$ws.Range("CF3").Value2 = 1509.0648000000001
$ws.Range("CG3").Value2 = 944.88
$ws.Range("CH3").Value2 = 1524
$ws.Range("CI3").Value2 = 1798.9295999999999
$ws.Range("CD4").Value2 = 767.38255033557004
$ws.Range("CE4").Value2 = 1250.2
$ws.Range("CF4").Value2 = 1481.1428571428601
$ws.Range("CG4").Value2 = 850
$ws.Range("CH4").Value2 = 1550
$ws.Range("CI4").Value2 = 1840
$ws.Range("CD5").Value2 = 365.43209876543199
$ws.Range("CE5").Value2 = 582.758620689655
$ws.Range("CF5").Value2 = 587.5
$ws.Range("CG5").Value2 = 300
$ws.Range("CH5").Value2 = 700
$ws.Range("CI5").Value2 = 650
$ws.Range("CD6").Value2 = 972.40460587388498
$ws.Range("CE6").Value2 = 1652.42651844286
$ws.Range("CF6").Value2 = 1968.86873612
$ws.Range("CG6").Value2 = 1066.799966
$ws.Range("CH6").Value2 = 1813.5599420999999
$ws.Range("CI6").Value2 = 3115.0559007000002
$ws.Range("CD8").Value2 = 247.11538461538501
$ws.Range("CE8").Value2 = 491.25
$ws.Range("CF8").Value2 = 925
$ws.Range("CG8").Value2 = 0
$ws.Range("CH8").Value2 = 537.5
$ws.Range("CI8").Value2 = 987.5
$ws.Range("CD9").Value2 = 353.41603053435102
$ws.Range("CE9").Value2 = 1243.2142857142901
$ws.Range("CF9").Value2 = 1608.75
$ws.Range("CG9").Value2 = 365
$ws.Range("CH9").Value2 = 1175
$ws.Range("CI9").Value2 = 1750
$ws.Range("CD10").Value2 = 302.564102564103
$ws.Range("CE10").Value2 = 373.33333333333297
$ws.Range("CF10").Value2 = 516.66666666666697
$ws.Range("CG10").Value2 = 300
$ws.Range("CH10").Value2 = 400
$ws.Range("CI10").Value2 = 550
$ws.Range("CD11").Value2 = 601.9
$ws.Range("CE11").Value2 = 755.83333333333303
$ws.Range("CF11").Value2 = 1300
$ws.Range("CG11").Value2 = 500
$ws.Range("CH11").Value2 = 475
$ws.Range("CI11").Value2 = 1350
$ws.Range("CD13").Value2 = 636.92307692307702
$ws.Range("CE13").Value2 = 852.94117647058795
$ws.Range("CF13").Value2 = 1025
$ws.Range("CG13").Value2 = 600
$ws.Range("CH13").Value2 = 1000
$ws.Range("CI13").Value2 = 1150
$ws.Range("CD14").Value2 = 499.97435897435901
$ws.Range("CE14").Value2 = 621.85
$ws.Range("CF14").Value2 = 791
$ws.Range("CG14").Value2 = 340
$ws.Range("CH14").Value2 = 562
$ws.Range("CI14").Value2 = 856.5
$ws.Range("CD15").Value2 = 378.57142857142901
$ws.Range("CE15").Value2 = 541.17647058823502
$ws.Range("CF15").Value2 = 950
$ws.Range("CG15").Value2 = 300
$ws.Range("CH15").Value2 = 400
$ws.Range("CI15").Value2 = 1350
$ws.Range("CD16").Value2 = 281.65680473372799
$ws.Range("CE16").Value2 = 477.19298245613999
$ws.Range("CF16").Value2 = 760
$ws.Range("CG16").Value2 = 200
$ws.Range("CH16").Value2 = 600
$ws.Range("CI16").Value2 = 800
$ws.Range("CD17").Value2 = 519.54022988505699
$ws.Range("CE17").Value2 = 725
$ws.Range("CF17").Value2 = 1040
$ws.Range("CG17").Value2 = 400
$ws.Range("CH17").Value2 = 800
$ws.Range("CI17").Value2 = 1100
$ws.Range("CD18").Value2 = 423.16708229426399
$ws.Range("CE18").Value2 = 890.57142857142901
$ws.Range("CF18").Value2 = 1367.5
$ws.Range("CG18").Value2 = 430
$ws.Range("CH18").Value2 = 1020
$ws.Range("CI18").Value2 = 1450
$ws.Range("CD21").Value2 = 580.12068965517199
$ws.Range("CE21").Value2 = 890.47619047619003
$ws.Range("CF21").Value2 = 1050.125
$ws.Range("CG21").Value2 = 650
$ws.Range("CH21").Value2 = 975
$ws.Range("CI21").Value2 = 1174.5

# --- Update view state: re-select the reported active cell in the
#     bottom-right (scrollable) pane so the saved selection matches the
#     target workbook. Freeze panes (xSplit=4 / ySplit=1) are already in
#     effect on this sheet and are left untouched.
$ws.Activate()
$ws.Range("CE4").Select()
